# New PO forecast model
# - Weekly Quantity: append 2 new weekly rows
# - Monthly Trend: append 1 new monthly row
# - PO Forecast: replace the forecasted quantities for all existing weeks,
#   shift the tail of the date series and append 2 new forecast weeks

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Weekly Quantity": add rows 75-76
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A75").Value = 45662.99999999999
$wsWeekly.Range("A75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B75").Value = 4
$wsWeekly.Range("A76").Value = 45676.99999999999
$wsWeekly.Range("A76").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B76").Value = 31

# ---------------------------------------------------------------------
# Sheet "Monthly Trend": add row 23
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A23").Value = 45688.99999999999
$wsMonthly.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMonthly.Range("B23").Value = 35

# ---------------------------------------------------------------------
# Sheet "PO Forecast": update forecast quantities for rows 2-74 (dates
# unchanged), overwrite the shifted date/quantity tail for rows 75-82,
# and append the two brand-new forecast weeks in rows 83-84.
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# New forecast quantities for rows 2 through 84 (row r -> index r-2)
$forecastQty = @(
    124,45,25,68,136,183,191,170,137,111,113,155,225,277,264,181,83,41,79,154,
    192,156,77,22,34,98,164,187,164,126,113,136,171,180,144,85,49,62,103,122,
    86,21,0,17,92,152,166,156,171,221,177,93,60,94,161,214,212,180,151,248,
    308,64,63,122,191,224,207,166,218,190,132,89,247,230,97,121,186,244,266,
    253,222,192,181
)

for ($i = 0; $i -lt $forecastQty.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).Value = $forecastQty[$i]
}

# New "ds" dates for the shifted/extended tail, rows 75-84
$forecastDates = @(
    45662.99999999999,
    45676.99999999999,
    45683.99999999999,
    45690.99999999999,
    45697.99999999999,
    45704.99999999999,
    45711.99999999999,
    45718.99999999999,
    45725.99999999999,
    45732.99999999999
)

for ($i = 0; $i -lt $forecastDates.Length; $i++) {
    $row = $i + 75
    $wsForecast.Cells.Item($row, 1).Value = $forecastDates[$i]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Output "PO forecast model updated"
